# Dropped in WRI's input data for first draft of 2.1. Made several fixes
# to model run errors.
#
# Target sheet: "CIRbTF" (Carbon Intensity Ratios by Transport Fuel).
# Adds three new fuel types (heavy/residual fuel oil, LPG propane or
# butane, hydrogen) with the same "-1" calculate-automatically flag as the
# existing fuels, clarifies the header label, and widens/wraps things so
# the longer header text fits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CIRbTF")

# Remember which sheet was active/selected so we can restore it - adding
# rows/selecting cells on CIRbTF below would otherwise leave CIRbTF as the
# active tab.
$originalActiveSheet = $wb.ActiveSheet.Name

# Add the new fuel rows *before* renaming the B1 header below, so the new
# fuel-name strings are registered ahead of the new header string (matches
# the order new entries were appended to the workbook's string table).
$ws.Range("A9").Value = "heavy or residual fuel oil"
$ws.Range("B9").Value = -1

$ws.Range("A10").Value = "LPG propane or butane"
$ws.Range("B10").Value = -1

$ws.Range("A11").Value = "hydrogen"
$ws.Range("B11").Value = -1

# Clarify the units of the header/flag column.
$ws.Range("B1").Value = "ratio or flag (dimensionless)"

# Let the now-longer header text wrap, and grow row 1 to fit it.
$ws.Range("B1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 45

# Column A needs to be a bit wider to comfortably fit the new, longer fuel
# names ("heavy or residual fuel oil", "LPG propane or butane").
$ws.Columns.Item(1).ColumnWidth = 25

# Leave the cursor on the header cell of the CIRbTF sheet.
$null = $ws.Range("B1").Select()

# Restore whichever sheet/tab was originally active.
$null = $wb.Worksheets.Item($originalActiveSheet).Activate()
